$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# IrrigationTask stack size increased: task1 stack formula G2 changes
# from $B$2+128 to $B$2+192 (dependent cells recalc automatically)
$ws.Range("G2").Formula = '=$B$2+192'

# Update the active selection to E16
$ws.Range("E16").Select()
